# This script reorders the 17 player data rows (A2:C18) on Sheet1.
# The header row (row 1) and the underlying set of (Player, Position, Team)
# triples are unchanged -- only the row order differs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Russell Westbrook", "PG", "Denver Nuggets"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Julian Champagnie", "SF,PF", "San Antonio Spurs"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Pascal Siakam", "SF,PF", "Indiana Pacers"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Bogdan Bogdanovic", "SG,SF", "Atlanta Hawks"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
